# Auto-generated edit script applying crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All D-column (Price) and E-column (Volume) values are stored as text
# in the source sheet, so force text number-format before assigning to
# prevent Excel from auto-converting numeric-looking strings into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.800.97"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.634.38"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.17"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.503"
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.257"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.83"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.643.83"
$ws.Range("E12").Value = "  +0.69%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.25"
$ws.Range("E13").Value = "  -0.59%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.860.90"
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("E15").Value = "  -0.40%  "
$ws.Range("E16").Value = "  +1.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.07"
$ws.Range("E17").Value = "  +0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.820.35"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "193.82"
$ws.Range("E20").Value = "  -0.53%  "
$ws.Range("E21").Value = "  +1.60%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.93"
$ws.Range("E22").Value = "  +0.83%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.14"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.77"
$ws.Range("E25").Value = "  -2.20%  "
$ws.Range("E26").Value = "  -0.16%  "
$ws.Range("E27").Value = "  -5.46%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.83"
$ws.Range("E29").Value = "  +0.45%  "
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0495"
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("E32").Value = "  +1.15%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +1.90%  "
$ws.Range("E35").Value = "  +0.11%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  +0.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.552"
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.111.45"
$ws.Range("E40").Value = "  +0.46%  "
$ws.Range("E41").Value = "  +0.73%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.57"
$ws.Range("E42").Value = "  +0.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.43"
$ws.Range("E43").Value = "  +2.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.795"
$ws.Range("E44").Value = "  -0.47%  "
$ws.Range("E45").Value = "  -2.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.40"
$ws.Range("E46").Value = "  +0.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.50"
$ws.Range("E47").Value = "  +12.26%  "
$ws.Range("E48").Value = "  -5.40%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.67"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0503"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("E51").Value = "  +0.28%  "
